$d = $word.ActiveDocument

# 1) "Ojo todo lo que..." -> "Ojo: (todo lo que..."
$d.Content.Find.Execute(
    "Ojo todo lo que está en parentesis no va en el código, son notas aclaratorias)",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Ojo: (todo lo que está en parentesis no va en el código, son notas aclaratorias)",
    2) | Out-Null

# 2) Fix typo "partuclarmente" -> "particularmente"
$d.Content.Find.Execute(
    "Puedes hacerlo de tres formas, (partuclarmente escogí la primera)",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Puedes hacerlo de tres formas, (particularmente escogí la primera)",
    2) | Out-Null

# 3) "Primera forma:" -> "Primera forma: (escogí esta)"
$d.Content.Find.Execute(
    "Primera forma:",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Primera forma: (escogí esta)",
    2) | Out-Null
